# Apply the crypto price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "26.003.78"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").Value = "1.744.74"
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.33"
$ws.Range("E5").Value = "  +2.26%  "

$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5053"
$ws.Range("E7").Value = "  -6.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2756"
$ws.Range("E8").Value = "  -2.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06195"
$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").Value = "1.744.84"
$ws.Range("E10").Value = "  -0.06%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07261"
$ws.Range("E11").Value = "  +0.96%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.20"
$ws.Range("E12").Value = "  -1.90%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6525"
$ws.Range("E13").Value = "  -0.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.697"
$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.63"
$ws.Range("E15").Value = "  -0.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9992"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "26.019.37"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.92"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006924"
$ws.Range("E20").Value = "  +1.51%  "

$ws.Range("D21").Value = "1.969.48"
$ws.Range("E21").Value = "  -0.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.485"
$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.764"
$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.370"
$ws.Range("E24").Value = "  +1.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.16"
$ws.Range("E25").Value = "  -2.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.507"
$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.783"
$ws.Range("E28").Value = "  -0.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.84"
$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.858"
$ws.Range("E30").Value = "  +0.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08198"
$ws.Range("E31").Value = "  -4.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.655"
$ws.Range("E32").Value = "  -1.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04667"
$ws.Range("E33").Value = "  +0.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.652"
$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9970"
$ws.Range("E35").Value = "  -0.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6108"
$ws.Range("E36").Value = "  -3.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.803"
$ws.Range("E37").Value = "  +3.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01625"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.935"
$ws.Range("E39").Value = "  -0.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9996"
$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.73"
$ws.Range("E41").Value = "  +0.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3930"
$ws.Range("E42").Value = "  -0.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7698"
$ws.Range("E43").Value = "  +2.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.009"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1166"
$ws.Range("E45").Value = "  +0.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.336"
$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.70"
$ws.Range("E47").Value = "  +0.83%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05324"
$ws.Range("E48").Value = "  -0.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.73"
$ws.Range("E49").Value = "  -1.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.627"
$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3470"
$ws.Range("E51").Value = "  -1.51%  "
